$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cos = $ws.ChartObjects()
$co = $cos.Item(2)
$chart = $co.Chart
$scol = $chart.SeriesCollection()
$s = $scol.Item(1)

$members = $s | Get-Member
$members | ForEach-Object { Write-Host $_.Name }
